$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") bumped from 46077 to 46078 for every data row (2-14)
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46078
}

# Rows 7-14 got their Beteckning/Datum/Area values shuffled between rows by the
# automatic source refresh. Apply the new row contents (col A, B, G) as per the
# updated data.
$rowData = @{
    7  = @("A 62884-2021", 44504, 0.8)
    8  = @("A 25634-2025", 45803.59570601852, 6)
    9  = @("A 28266-2025", 45818.56381944445, 1.9)
    10 = @("A 25015-2023", 45085.6989699074, 1.8)
    11 = @("A 19922-2025", 45771.63034722222, 10.1)
    12 = @("A 60024-2025", 45992, 1.1)
    14 = @("A 14271-2021", 44278, 6.7)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 7).Value = $vals[2]
}
